# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds header "K"; update the computed K values for rows 2-15.
$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 1
    6  = 0
    7  = 3
    8  = 1
    9  = 0
    10 = 6
    11 = 1
    12 = 1
    13 = 2
    14 = 2
    15 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
